# Generate Report for Archive
# "Ready for handoff" status updates to "In Translation" across the
# Overview, zh-cn and de-de sheets once the handback report has been
# (re-)generated for archival purposes.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5:F7").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5:C7").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5:C7").Value = "In Translation"
